$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '59.022.81'
$ws.Cells.Item(2, 5).Value = '  +0.28%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '2.500.44'
$ws.Cells.Item(3, 5).Value = '  -0.46%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.25%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '539.85'
$ws.Cells.Item(5, 5).Value = '  +1.39%  '

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '138.17'
$ws.Cells.Item(6, 5).Value = '  -0.29%  '

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.997'
$ws.Cells.Item(7, 5).Value = '  -0.29%  '

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.567'
$ws.Cells.Item(8, 5).Value = '  +1.00%  '

# Row 9
$ws.Cells.Item(9, 4).Value = '2.524.86'
$ws.Cells.Item(9, 5).Value = '  +0.38%  '

# Row 10
$ws.Cells.Item(10, 5).Value = '  +1.72%  '

# Row 11
$ws.Cells.Item(11, 5).Value = '  -0.79%  '

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '5.36'
$ws.Cells.Item(12, 5).Value = '  -1.45%  '

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '0.349'
$ws.Cells.Item(13, 5).Value = '  -2.02%  '

# Row 14
$ws.Cells.Item(14, 4).Value = '2.951.85'
$ws.Cells.Item(14, 5).Value = '  -0.16%  '

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '23.23'
$ws.Cells.Item(15, 5).Value = '  +0.76%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '58.949.86'
$ws.Cells.Item(16, 5).Value = '  +0.24%  '

# Row 17
$ws.Cells.Item(17, 5).Value = '  +0.08%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '2.527.30'
$ws.Cells.Item(18, 5).Value = '  +0.70%  '

# Row 19
$ws.Cells.Item(19, 5).Value = '  +0.91%  '

# Row 20
$ws.Cells.Item(20, 5).Value = '  +0.88%  '

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '325.90'
$ws.Cells.Item(21, 5).Value = '  +1.06%  '

# Row 22
$ws.Cells.Item(22, 5).Value = '  +0.04%  '

# Row 23
$ws.Cells.Item(23, 5).Value = '  +2.80%  '

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '65.37'
$ws.Cells.Item(24, 5).Value = '  +5.22%  '

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '0.424'
$ws.Cells.Item(25, 5).Value = '  -0.14%  '

# Row 26
$ws.Cells.Item(26, 5).Value = '  +0.70%  '

# Row 27
$ws.Cells.Item(27, 5).Value = '  +0.51%  '

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '7.69'
$ws.Cells.Item(28, 5).Value = '  -0.89%  '

# Row 29
$ws.Cells.Item(29, 4).Value = '0.0₃0780'
$ws.Cells.Item(29, 5).Value = '  +1.47%  '

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '6.72'
$ws.Cells.Item(30, 5).Value = '  +0.77%  '

# Row 31
$ws.Cells.Item(31, 5).Value = '  +0.45%  '

# Row 32
$ws.Cells.Item(32, 2).Value = 'Fetch.AI'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '1.19'
$ws.Cells.Item(32, 5).Value = '  +6.87%  '

# Row 33
$ws.Cells.Item(33, 2).Value = 'Monero'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '166.64'
$ws.Cells.Item(33, 5).Value = '  +1.92%  '

# Row 34
$ws.Cells.Item(34, 2).Value = 'USDe'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '0.999'
$ws.Cells.Item(34, 5).Value = '  +0.00%  '

# Row 35
$ws.Cells.Item(35, 2).Value = 'ImmutableX'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '1.47'
$ws.Cells.Item(35, 5).Value = '  +3.41%  '

# Row 36
$ws.Cells.Item(36, 5).Value = '  +0.43%  '

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '4.14'
$ws.Cells.Item(37, 5).Value = '  -2.21%  '

# Row 38
$ws.Cells.Item(38, 5).Value = '  -0.17%  '

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '36.86'
$ws.Cells.Item(39, 5).Value = '  +0.16%  '

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '0.823'
$ws.Cells.Item(40, 5).Value = '  +2.57%  '

# Row 41
$ws.Cells.Item(41, 5).Value = '  +0.23%  '

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '288.50'
$ws.Cells.Item(42, 5).Value = '  +3.43%  '

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '5.27'
$ws.Cells.Item(43, 5).Value = '  +1.02%  '

# Row 44
$ws.Cells.Item(44, 2).Value = 'Aave'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '131.80'
$ws.Cells.Item(44, 5).Value = '  +8.34%  '

# Row 45
$ws.Cells.Item(45, 2).Value = 'Mantle'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '0.610'
$ws.Cells.Item(45, 5).Value = '  +2.43%  '

# Row 46
$ws.Cells.Item(46, 2).Value = 'FirstDigitalUSD'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '0.994'
$ws.Cells.Item(46, 5).Value = '  -0.49%  '

# Row 47
$ws.Cells.Item(47, 5).Value = '  +0.01%  '

# Row 48
$ws.Cells.Item(48, 5).Value = '  +0.11%  '

# Row 49
$ws.Cells.Item(49, 5).Value = '  +0.39%  '

# Row 50
$ws.Cells.Item(50, 5).Value = '  +0.00%  '

# Row 51
$ws.Cells.Item(51, 5).Value = '  -0.92%  '
